$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Make room in the data table: insert two fresh rows right above the first
#    data row (16). This shifts the existing 10 data rows (old 16-25,
#    VICTOR MANUEL IBARRA MIRANDA / periods 2303-2312) down to rows 18-27,
#    keeping their original styling/merges intact (the "last row" bottom
#    border naturally rides along on the former row 25 -> now row 27), and
#    also pushes the footer (signature block) from rows 30/31 down to 32/33
#    as a side effect (it sits below the inserted rows).
# ---------------------------------------------------------------------------
$ws.Rows.Item(16).Insert()
$ws.Rows.Item(16).Insert()

# The two newly inserted blank rows (16 & 17) pick up a generic style; copy
# the exact formatting used by the rest of the data rows (row 18, which is
# the "middle" row style) onto them so every data row matches.
$ws.Range("B18:J18").Copy() | Out-Null
$ws.Range("B16:J17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3) New employee: two rows (1047477096 / MARLYN ISABEL CORRALES AMOR),
#    periods 2504 and 2503.
# ---------------------------------------------------------------------------
$ws.Cells.Item(16, 2).Value2 = "CC"
$ws.Cells.Item(16, 3).Value2 = "1047477096"
$ws.Cells.Item(16, 4).Value2 = "MARLYN ISABEL CORRALES AMOR"
$ws.Cells.Item(16, 5).Value2 = "2504"
$ws.Cells.Item(16, 6).Value2 = 56940
$ws.Cells.Item(16, 7).Value2 = 1300000

$ws.Cells.Item(17, 2).Value2 = "CC"
$ws.Cells.Item(17, 3).Value2 = "1047477096"
$ws.Cells.Item(17, 4).Value2 = "MARLYN ISABEL CORRALES AMOR"
$ws.Cells.Item(17, 5).Value2 = "2503"
$ws.Cells.Item(17, 6).Value2 = 56940
$ws.Cells.Item(17, 7).Value2 = 1300000

# ---------------------------------------------------------------------------
# 4) Existing employee (19895150 / VICTOR MANUEL IBARRA MIRANDA): periods now
#    listed in descending order, 2312 down to 2303, across rows 18-27.
# ---------------------------------------------------------------------------
$periods = @(2312, 2311, 2310, 2309, 2308, 2307, 2306, 2305, 2304, 2303)
for ($i = 0; $i -lt $periods.Length; $i++) {
    $r = 18 + $i
    $ws.Cells.Item($r, 2).Value2 = "CC"
    $ws.Cells.Item($r, 3).Value2 = "19895150"
    $ws.Cells.Item($r, 4).Value2 = "VICTOR MANUEL IBARRA MIRANDA"
    $ws.Cells.Item($r, 5).Value2 = [string]$periods[$i]
    $ws.Cells.Item($r, 7).Value2 = 1160000
}
# All periods keep the 46400 "Valor Mora" value, except the last one (period
# 2303, now row 27) which keeps the original 43307 value.
for ($r = 18; $r -le 26; $r++) {
    $ws.Cells.Item($r, 6).Value2 = 46400
}
$ws.Cells.Item(27, 6).Value2 = 43307

# ---------------------------------------------------------------------------
# 5) Header / summary fields.
# ---------------------------------------------------------------------------
$ws.Cells.Item(11, 5).Value2 = 574787   # VALOR MORA total
$ws.Cells.Item(13, 3).Value2 = 2        # Cant. Trabajadores
$ws.Cells.Item(13, 6).Value2 = 12       # Cant. Periodos

Write-Output "edit complete"
